# "finished logins, weibo share"
#
# Column B on Sheet1 holds each module's completion percentage; column D
# holds a free-text "remaining work" note.
# - D14 (note for the "登录/注册" / login-register row) text updated: the
#   generic placeholder "待接通接口" is replaced with the concrete
#   remaining task "待接通第三方登陆，KidsTC注册"
# - B14 (progress for "登录/注册") bumped to 100% (finished logins) and
#   gets the same "done" highlight fill used by the other 100%-complete rows
# - B13 (progress for "支付") nudged from blank to 10%
# - B18 (progress for "分享"/share, i.e. weibo share) bumped from 20% to 60%
# - selection cursor ends up parked on D20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Login/register row: finished -> update remaining-work note and progress
$ws.Range("D14").Value = "待接通第三方登陆，KidsTC注册"
$ws.Range("B14").Value = 1

# Copy the "complete" cell formatting (green fill) from an already-100%
# cell onto B14 so it matches the rest of the finished rows.
$ws.Range("B3").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# Payment row note progress nudged up
$ws.Range("B13").Value = 0.1

# Share (weibo share) row progress nudged up
$ws.Range("B18").Value = 0.6

# Leave the selection where the author last left it
$ws.Range("D20").Select()
